$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet is protected; unprotect before editing
$ws.Unprotect("D382")

# Update the confidential/model-holdings date string in A13
$ws.Range("A13").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-14 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) values for rows 2-10
$ws.Range("D2").Value = 0.08808259255053792
$ws.Range("E2").Value = 0.03965236284627927

$ws.Range("D3").Value = 0.1052292869178707
$ws.Range("E3").Value = 0.02347652347652329

$ws.Range("D4").Value = 0.1211859218529886
$ws.Range("E4").Value = 0.01402535657686199

$ws.Range("D5").Value = 0.1420294826105429
$ws.Range("E5").Value = 0.01718019035328289

$ws.Range("D6").Value = 0.1395691728783703
$ws.Range("E6").Value = 0.008328706274292053

$ws.Range("D7").Value = 0.1482917828956922
$ws.Range("E7").Value = 0.01488691669052389

$ws.Range("D8").Value = 0.1256975005045623
$ws.Range("E8").Value = 0.02527743526510484

$ws.Range("D9").Value = 0.1299142597894349
$ws.Range("E9").Value = 0.02101030927835046

$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 0.01937975733165209

# Re-protect the sheet as it was originally
$ws.Protect("D382", $false, $true, $true, $true, $false, $false)
